$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D (shifting the old D -> E) so the old D data (sums to 100)
# moves into the new E column, and new values are put into D.
$ws.Columns("D").Insert()

# Header for the new D column - actually wait, we need new E header "g4"
# Let's set the header in E1 first
$ws.Range("E1").Value = "g4"

# Now set new values for column D (rows 2-8)
$ws.Range("D2").Value = 15
$ws.Range("D3").Value = 10
$ws.Range("D4").Value = 20
$ws.Range("D5").Value = 10
$ws.Range("D6").Value = 10
$ws.Range("D7").Value = 20
$ws.Range("D8").Value = 15

# Formula for D9 and E9 sums
$ws.Range("D9").Formula = "=SUM(D2:D8)"
$ws.Range("E9").Formula = "=SUM(E2:E8)"

# Conditional formatting: copy format from D to E (since insert should have copied it to E already)
# But we need D9 to have its own separate conditional formatting rule too.
# Let's check existing conditional formatting and split.
